$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column for rows 2-29 from 45581 to 45582
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 3).Value = 45582
}

# Remove the last data row (row 30, "A 45908-2024"), shifting rows up
$ws.Rows(30).Delete()
